$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (sharedStrings rich-text runs) ---
$ws.Range("A8").Value = "Volume 32   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# --- Cells changing type/style (copy format from a stable donor cell, then set value) ---
$ws.Range("F15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 2

$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("F15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1

$ws.Range("M14").Copy()
$ws.Range("M22").PasteSpecial(-4122)
$ws.Range("M22").Value = 200

$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("F15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 2

$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1

$ws.Range("F15").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1

$ws.Range("M14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 0

$ws.Range("F15").Copy()
$ws.Range("J28").PasteSpecial(-4122)
$ws.Range("J28").Value = 1

$ws.Range("M14").Copy()
$ws.Range("K28").PasteSpecial(-4122)
$ws.Range("K28").Value = 100

$ws.Range("F15").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("C31").Value = 1

$ws.Range("F15").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F31").Value = 1

$ws.Range("C14").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$ws.Range("C14").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$ws.Range("F15").Copy()
$ws.Range("I31").PasteSpecial(-4122)
$ws.Range("I31").Value = 1

# --- Cells with same style (direct value assignment) ---
$ws.Range("I15").Value = 4
$ws.Range("K15").Value = 300
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -11.111111111111
$ws.Range("I16").Value = 11
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = -8.333333333333
$ws.Range("L16").Value = -31.25
$ws.Range("M16").Value = -15.384615384615
$ws.Range("N16").Value = -81.666666666666
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -41.176470588235
$ws.Range("I17").Value = 13
$ws.Range("J17").Value = 18
$ws.Range("K17").Value = -27.777777777777
$ws.Range("L17").Value = -23.529411764705
$ws.Range("M17").Value = 62.5
$ws.Range("N17").Value = -35
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 4
$ws.Range("K18").Value = 50
$ws.Range("L18").Value = -25
$ws.Range("M18").Value = -33.333333333333
$ws.Range("N18").Value = -91.304347826087
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = -38.461538461538
$ws.Range("I19").Value = 21
$ws.Range("J19").Value = 29
$ws.Range("K19").Value = -27.586206896551
$ws.Range("L19").Value = -27.586206896551
$ws.Range("M19").Value = -4.545454545454
$ws.Range("N19").Value = -68.181818181818
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -90
$ws.Range("N20").Value = -97.727272727272
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -12.5
$ws.Range("F21").Value = 42
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = -27.586206896551
$ws.Range("I21").Value = 56
$ws.Range("J21").Value = 66
$ws.Range("K21").Value = -15.151515151515
$ws.Range("L21").Value = -30
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = -78.867924528301
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 3
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 400
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 20
$ws.Range("I23").Value = 14
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = 7.692307692307
$ws.Range("L23").Value = -12.5
$ws.Range("M23").Value = 75
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 60
$ws.Range("I24").Value = 35
$ws.Range("J24").Value = 27
$ws.Range("K24").Value = 29.629629629629
$ws.Range("L24").Value = -12.5
$ws.Range("M24").Value = -20.454545454545
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 20
$ws.Range("I25").Value = 7
$ws.Range("J25").Value = 7
$ws.Range("L25").Value = -61.111111111111
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 19
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = -40.625
$ws.Range("J26").Value = 34
$ws.Range("K26").Value = -26.470588235294
$ws.Range("L26").Value = 25
$ws.Range("M26").Value = 38.888888888888
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 5
$ws.Range("K27").Value = 150
$ws.Range("L27").Value = 400
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 2
$ws.Range("L28").Value = -66.666666666666

$excel.CutCopyMode = $false
